$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text storage (avoid Excel auto-converting numeric-looking strings)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.023.77"
$ws.Range("E2").Value = "  +1.45%  "

Set-TextValue $ws.Range("D3") "2.588.57"
$ws.Range("E3").Value = "  -0.19%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.21%  "

Set-TextValue $ws.Range("D5") "527.28"
$ws.Range("E5").Value = "  +0.98%  "

Set-TextValue $ws.Range("D6") "139.18"
$ws.Range("E6").Value = "  -2.73%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.71%  "

Set-TextValue $ws.Range("D9") "2.601.01"
$ws.Range("E9").Value = "  -0.51%  "

Set-TextValue $ws.Range("D10") "6.50"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  -0.21%  "

Set-TextValue $ws.Range("D12") "0.331"
$ws.Range("E12").Value = "  -3.07%  "

$ws.Range("E13").Value = "  +3.08%  "

Set-TextValue $ws.Range("D14") "3.049.55"
$ws.Range("E14").Value = "  -0.18%  "

Set-TextValue $ws.Range("D15") "58.958.35"
$ws.Range("E15").Value = "  +1.37%  "

Set-TextValue $ws.Range("D16") "20.52"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "2.599.13"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.0000133"
$ws.Range("E18").Value = "  -0.87%  "

Set-TextValue $ws.Range("D19") "342.97"
$ws.Range("E19").Value = "  +1.00%  "

Set-TextValue $ws.Range("D20") "4.31"
$ws.Range("E20").Value = "  -0.75%  "

Set-TextValue $ws.Range("D21") "10.08"
$ws.Range("E21").Value = "  -1.53%  "

Set-TextValue $ws.Range("D22") "6.43"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +0.10%  "

Set-TextValue $ws.Range("D24") "66.36"
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("E29").Value = "  +0.03%  "

Set-TextValue $ws.Range("D30") "0.0₃0722"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.61"
$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D32") "5.91"
$ws.Range("E32").Value = "  -4.18%  "

Set-TextValue $ws.Range("D33") "18.72"
$ws.Range("E33").Value = "  -0.35%  "

Set-TextValue $ws.Range("D34") "149.36"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("E35").Value = "  -0.97%  "

$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D37") "36.82"
$ws.Range("E37").Value = "  +2.35%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D38") "1.49"
$ws.Range("E38").Value = "  +1.76%  "

Set-TextValue $ws.Range("D39") "0.825"
$ws.Range("E39").Value = "  -4.53%  "

Set-TextValue $ws.Range("D40") "0.810"
$ws.Range("E40").Value = "  -6.46%  "

Set-TextValue $ws.Range("D41") "3.52"
$ws.Range("E41").Value = "  -0.50%  "

Set-TextValue $ws.Range("D42") "0.998"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.602"
$ws.Range("E43").Value = "  -0.67%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D44") "270.39"
$ws.Range("E44").Value = "  -0.14%  "

Set-TextValue $ws.Range("D45") "10.73"
$ws.Range("E45").Value = "  +0.59%  "

Set-TextValue $ws.Range("D46") "0.0955"
$ws.Range("E46").Value = "  -0.20%  "

Set-TextValue $ws.Range("D47") "0.0514"
$ws.Range("E47").Value = "  -1.46%  "

Set-TextValue $ws.Range("D48") "18.40"
$ws.Range("E48").Value = "  -1.79%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "4.63"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D50") "1.964.04"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("E51").Value = "  +0.11%  "
